$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- 1. Update the lookup table on Sheet3 (A20:B36) with the new values.
#     Sheet3!C2:C18 (VLOOKUP against this table) and, in turn,
#     Sheet1!CB2:CC18 (VLOOKUP against Sheet3!B1:C18) recalculate
#     automatically because they are formulas.
$newLookup = @{
    20 = 10.598440240581542
    21 = 0.52491621463652338
    22 = 7.7689023041750893
    23 = 3.7025131645317715
    24 = 7.4118289328797315
    25 = 0
    26 = 6.1002054852828183
    27 = 15.352112605727879
    28 = 5.3076711488414574
    29 = 3.0727409418011855
    30 = 10.304569300734874
    31 = 18.938279506344095
    32 = 4.5500612098362145
    33 = 8.7175212733279697
    34 = 6.8722895180192314
    35 = 6.3752036495534714
    36 = 51.15410626449934
}

foreach ($r in $newLookup.Keys) {
    $ws3.Cells.Item($r, 2).Value = $newLookup[$r]
}

$wb.Application.Calculate()

# --- 2. Add a new "06-nov" snapshot column (CI) on Sheet1, mirroring the
#     existing CH ("03-nov") column, but holding today's freshly
#     recalculated VLOOKUP values (CB/CC) as static numbers.
$ws1.Cells.Item(1, 87).Value = "06-nov"
$ws1.Cells.Item(1, 87).NumberFormat = "@"

for ($r = 2; $r -le 18; $r++) {
    $ws1.Cells.Item($r, 87).Value = $ws1.Cells.Item($r, 80).Value2
    $ws1.Cells.Item($r, 87).NumberFormat = "0"
}

# --- 3. Move the active selection, matching the post-edit workbook state.
$ws1.Range("CB2").Select()
